# Bug fix for WinPt Check:
#  - Column A (the WinPt index) had been left blank for rows 46-57 because the
#    sheet ran out of WinPt's before the SGConfig list did. Backfill those
#    cells with the expected sequential index values (43-54) so the lookup
#    logic has something to compare against instead of treating them as
#    undefined.
#  - Refresh the view so it reflects where the fix was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstRow = 46
$firstValue = 43
$lastRow = 57

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $value = $firstValue + ($row - $firstRow)
    $ws.Cells.Item($row, 1).Value = $value
}

# Update the active view: scrolled position and current selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B3").Select()
